$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 42 has the same even-row formatting (shaded fill + borders) that the
# new row should inherit, so copy its formatting down to row 44 first.
$ws.Range("A42:F42").Copy()
$ws.Range("A44").PasteSpecial(-4122)

$ws.Range("A44").Value = "F"
$ws.Range("B44").Value = "UREQ026"
$ws.Range("C44").Value = ""
$ws.Range("D44").Value = "There will be 10 sets of 'T'-shaped words where the first letter of the vertical word overlaps with the horizontal word"
$ws.Range("E44").Value = "Product Description"
$ws.Range("F44").Value = "WordBank TWords()"

$ws.Range("A44:F44").Select()
